# Regenerate the localization-status report values/column widths
# (mirrors what the "Generate Report for Archive" report generator does
#  when it re-runs: it refreshes the Status text and re-sizes the
#  Status-related columns to fit the new text).

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$zhcn     = $wb.Sheets.Item("zh-cn")
$dede     = $wb.Sheets.Item("de-de")

# The "Status" value shared across the Overview roll-up (columns E/F)
# and each language sheet's Status column (column C) moves from
# "Ready for handoff" to "In Translation".
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Re-fit the Status columns now that the text is shorter.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
